# Login.xlsx - "Registor" sheet data + formatting update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix verificationMessage typo ("Welcome123" -> "Welcome") in rows 3 & 5 ---
$msg = "Welcome to your account. Here you can manage all of your personal information and orders."
$ws.Range("J3").Value = $msg
$ws.Range("J5").Value = $msg

# --- Update selectYear (column D) data values for rows 2-5 ---
$ws.Range("D2").Value = 2017
$ws.Range("D3").Value = 2016
$ws.Range("D4").Value = 2015
$ws.Range("D5").Value = 2014

# --- Highlight the updated D and F columns (rows 2-5) with a yellow fill ---
$ws.Range("D2:D5").Interior.Color = 65535
$ws.Range("F2:F5").Interior.Color = 65535

# --- Move/expand the active selection to F2:F5 ---
$ws.Activate()
[void]$ws.Range("F2:F5").Select()
